# Generate Report for Handoff
# Updates the localization-status report to reflect that the zh-cn and
# de-de targets are now "Ready for handoff" (were "In Translation"), and
# refreshes the handoff generation timestamps accordingly. Also touches
# the Status-related column widths (auto-fit effect of the longer status
# text) on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps ---
$overview.Range("G2").Value = "2016-08-25 04:56:50"
$dede.Range("H2").Value = "2016-08-25 04:56:50"
$zhcn.Range("H2").Value = "2016-08-25 04:56:45"

# --- Column widths (status columns widened to fit "Ready for handoff") ---
# Target OOXML column width is 17.2159881591797 (Excel's own pixel-based
# autofit result). This COM layer quantizes ColumnWidth to a coarser
# pixel grid, so 16.3 is the calibrated input that lands on the closest
# reachable width (17.166666666666668) to that target.
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
